$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes -------------------------------------------------
# The "2 minutes running" audio prompt now also mentions incline (độ dốc),
# not just speed (tốc độ).
$ws.Range("C20").Value = "Bạn đã chạy được 2 phút với tốc độ và độ dốc hiện tại, hãy bấm phím speed cộng để tăng hoặc speed trừ để giảm tốc độ, bấm phím incline cộng để tăng hoặc incline trừ để giảm độ dốc"
$ws.Rows(20).RowHeight = 30

# New audio prompt: incline increased warning (row 25 gains a Content cell)
$ws.Range("C25").Value = "Tăng độ dốc, chú ý giữ an toàn nhé"

# New audio prompt: incline decreased warning (row 26 becomes a full data row)
$ws.Range("A26").Value = 22
$ws.Range("B26").Value = "025"
$ws.Range("C26").Value = "Giảm độ dốc, chú ý giữ an toàn nhé"

# --- Row-height touch-ups (match re-measured wrap heights) -----------
$ws.Rows(2).RowHeight = 90
$ws.Rows(3).RowHeight = 30
$ws.Rows(4).RowHeight = 30
$ws.Rows(5).RowHeight = 30
$ws.Rows(6).RowHeight = 30
$ws.Rows(7).RowHeight = 30
$ws.Rows(8).RowHeight = 30
$ws.Rows(9).RowHeight = 45
$ws.Rows(10).RowHeight = 60
$ws.Rows(11).RowHeight = 45

# --- View state: selection follows the last edited content row -------
$ws.Range("C26").Select()
